$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 11:20"

# Full refreshed country data (rows 4-205): País, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$data = @(
    @(4, "Estados Unidos", 123780, 202, 3238, 118313, 2666, 8, 2229),
    @(5, "Italia", 92472, 0, 12384, 70065, 3856, 0, 10023),
    @(6, "China", 81439, 45, 75448, 2691, 742, 5, 3300),
    @(7, "España", 73235, 0, 12285, 54968, 4165, 0, 5982),
    @(8, "Alemania", 58247, 552, 8481, 49311, 1581, 22, 455),
    @(9, "Francia", 37575, 0, 5700, 29561, 4273, 0, 2314),
    @(10, "Iran", 35408, 0, 11679, 21212, 3206, 0, 2517),
    @(11, "Reino Unido", 17089, 0, 135, 15935, 163, 0, 1019),
    @(12, "Suiza", 14076, 0, 1595, 12217, 301, 0, 264),
    @(13, "Belgica", 10836, 1702, 1359, 9046, 867, 78, 431),
    @(14, "Paises Bajos", 9762, 0, 3, 9120, 914, 0, 639),
    @(15, "Corea del Sur", 9583, 105, 5033, 4398, 59, 8, 152),
    @(16, "Austria", 8395, 124, 225, 8102, 135, 0, 68),
    @(17, "Turquia", 7402, 0, 70, 7224, 445, 0, 108),
    @(18, "Canada", 5655, 0, 508, 5087, 120, 0, 60),
    @(19, "Portugal", 5170, 0, 43, 5027, 89, 0, 100),
    @(20, "Noruega", 4043, 28, 7, 4013, 84, 0, 23),
    @(21, "Australia", 3969, 334, 226, 3727, 23, 2, 16),
    @(22, "Brasil", 3904, 0, 6, 3784, 296, 0, 114),
    @(23, "Israel", 3865, 246, 89, 3764, 66, 0, 12),
    @(24, "Suecia", 3447, 0, 16, 3326, 239, 0, 105),
    @(25, "Chequia", 2669, 38, 11, 2645, 45, 2, 13),
    @(26, "Malasia", 2470, 150, 320, 2123, 73, 0, 27),
    @(27, "Irlanda", 2415, 0, 5, 2374, 59, 0, 36),
    @(28, "Dinamarca", 2201, 0, 1, 2135, 109, 0, 65),
    @(29, "Chile", 1909, 0, 61, 1842, 7, 0, 6),
    @(30, "Luxemburgo", 1831, 0, 40, 1773, 25, 0, 18),
    @(31, "Ecuador", 1823, 0, 3, 1772, 58, 0, 48),
    @(32, "Polonia", 1717, 79, 7, 1691, 3, 1, 19),
    @(33, "Japon", 1693, 0, 404, 1237, 56, 0, 52),
    @(34, "Rusia", 1534, 270, 64, 1462, 8, 4, 8),
    @(35, "Pakistan", 1526, 31, 29, 1484, 11, 1, 13),
    @(36, "Rumania", 1452, 0, 139, 1276, 34, 0, 37),
    @(37, "Filipinas", 1418, 343, 42, 1305, 1, 3, 71),
    @(38, "Tailandia", 1388, 143, 97, 1284, 11, 1, 7),
    @(39, "Indonesia", 1285, 130, 64, 1107, 0, 12, 114),
    @(40, "Finlandia", 1221, 54, 10, 1202, 32, 0, 9),
    @(41, "Arabia Saudita", 1203, 0, 37, 1162, 6, 0, 4),
    @(42, "Sudafrica", 1187, 0, 31, 1155, 7, 0, 1),
    @(43, "Grecia", 1061, 0, 52, 977, 66, 0, 32),
    @(44, "India", 987, 0, 87, 875, 0, 1, 25),
    @(45, "Islandia", 963, 0, 114, 847, 18, 0, 2),
    @(46, "Panama", 901, 0, 4, 880, 32, 0, 17),
    @(47, "Mexico", 848, 131, 4, 828, 1, 4, 16),
    @(48, "Singapur", 802, 0, 198, 601, 19, 1, 3),
    @(49, "Argentina", 745, 0, 72, 654, 0, 0, 19),
    @(50, "Republica Dominicana", 719, 0, 3, 688, 0, 0, 28),
    @(51, "Crucero", 712, 0, 597, 105, 15, 0, 10),
    @(52, "Eslovenia", 684, 0, 10, 665, 25, 0, 9),
    @(53, "Peru", 671, 0, 16, 639, 33, 0, 16),
    @(54, "Serbia", 659, 0, 42, 606, 25, 1, 11),
    @(55, "Croacia", 657, 0, 45, 607, 14, 0, 5),
    @(56, "Estonia", 645, 0, 20, 624, 10, 0, 1),
    @(57, "Colombia", 608, 0, 10, 592, 0, 0, 6),
    @(58, "Catar", 590, 0, 45, 544, 6, 0, 1),
    @(59, "Hong Kong", 582, 22, 112, 466, 5, 0, 4),
    @(60, "Egipto", 576, 0, 121, 419, 0, 0, 36),
    @(61, "Nueva Zelanda", 514, 0, 56, 457, 1, 0, 1),
    @(62, "Irak", 506, 0, 131, 333, 0, 0, 42),
    @(63, "Barein", 476, 0, 265, 207, 1, 0, 4),
    @(64, "Emiratos Arabes Unidos", 468, 0, 55, 411, 2, 0, 2),
    @(65, "Argelia", 454, 0, 31, 394, 0, 0, 29),
    @(66, "Lituania", 437, 43, 1, 429, 2, 0, 7),
    @(67, "Marruecos", 437, 35, 12, 399, 1, 1, 26),
    @(68, "Armenia", 424, 17, 30, 391, 6, 2, 3),
    @(69, "Ucrania", 418, 62, 5, 404, 0, 0, 9),
    @(70, "Libano", 412, 0, 30, 374, 4, 0, 8),
    @(71, "Hungria", 408, 65, 34, 361, 6, 2, 13),
    @(72, "Letonia", 347, 42, 1, 346, 3, 0, 0),
    @(73, "Bulgaria", 338, 7, 11, 319, 8, 1, 8),
    @(74, "Principado de Andorra", 308, 0, 1, 304, 10, 0, 3),
    @(75, "Uruguay", 304, 0, 0, 303, 9, 0, 1),
    @(76, "Taiwan", 298, 15, 30, 266, 0, 0, 2),
    @(77, "Costa Rica", 295, 0, 3, 290, 6, 0, 2),
    @(78, "Eslovaquia", 292, 0, 2, 290, 1, 0, 0),
    @(79, "Tunez", 278, 0, 2, 268, 10, 0, 8),
    @(80, "Bosnia y Herzegovina", 278, 0, 8, 264, 1, 0, 6),
    @(81, "Kazajistan", 251, 23, 18, 232, 0, 0, 1),
    @(82, "Jordania", 246, 0, 18, 227, 3, 0, 1),
    @(83, "Republica de Macedonia", 241, 0, 3, 234, 1, 0, 4),
    @(84, "Kuwait", 235, 0, 67, 168, 12, 0, 0),
    @(85, "Moldavia", 231, 0, 2, 227, 33, 0, 2),
    @(86, "San Marino", 224, 0, 6, 196, 16, 0, 22),
    @(87, "Burkina Faso", 207, 0, 21, 175, 0, 0, 11),
    @(88, "Albania", 197, 0, 33, 154, 3, 0, 10),
    @(89, "Reunion", 183, 0, 1, 182, 0, 0, 0),
    @(90, "Azerbaiyan", 182, 0, 15, 163, 23, 0, 4),
    @(91, "Republica de Chipre", 179, 0, 15, 159, 3, 0, 5),
    @(92, "Vietnam", 179, 5, 21, 158, 3, 0, 0),
    @(93, "Oman", 167, 15, 23, 144, 0, 0, 0),
    @(94, "Islas Feroe", 159, 4, 70, 89, 1, 0, 0),
    @(95, "Malta", 149, 0, 2, 147, 1, 0, 0),
    @(96, "Ghana", 141, 0, 2, 134, 1, 0, 5),
    @(97, "Uzbekistan", 133, 29, 7, 124, 8, 0, 2),
    @(98, "Senegal", 130, 0, 18, 112, 0, 0, 0),
    @(99, "Brunei", 120, 0, 25, 94, 1, 0, 1),
    @(100, "Cuba", 119, 0, 4, 112, 2, 0, 3),
    @(101, "Venezuela", 119, 0, 39, 78, 2, 0, 2),
    @(102, "Sri Lanka", 115, 2, 10, 104, 5, 0, 1),
    @(103, "Honduras", 110, 15, 3, 106, 4, 0, 1),
    @(104, "Afganistan", 110, 0, 2, 104, 0, 0, 4),
    @(105, "Estado de Palestina", 104, 0, 18, 85, 0, 0, 1),
    @(106, "Camboya", 103, 4, 21, 82, 1, 0, 0),
    @(107, "Mauricio", 102, 0, 0, 100, 1, 0, 2),
    @(108, "Guadalupe", 102, 0, 17, 83, 4, 0, 2),
    @(109, "Costa de Marfil", 101, 0, 3, 98, 0, 0, 0),
    @(110, "Nigeria", 97, 0, 3, 93, 0, 0, 1),
    @(111, "Bielorrusia", 94, 0, 32, 62, 2, 0, 0),
    @(112, "Martinica", 93, 0, 0, 92, 12, 0, 1),
    @(113, "Camerun", 91, 0, 2, 87, 0, 0, 2),
    @(114, "Georgia", 90, 0, 18, 72, 1, 0, 0),
    @(115, "Kirguistan", 84, 26, 0, 84, 0, 0, 0),
    @(116, "Montenegro", 84, 0, 0, 83, 1, 0, 1),
    @(117, "Bolivia", 81, 7, 0, 81, 3, 0, 0),
    @(118, "Trinidad yTobago", 76, 0, 1, 72, 0, 0, 3),
    @(119, "Consejo Danes para los Refugiados", 65, 0, 2, 57, 0, 0, 6),
    @(120, "Mayotte", 63, 0, 0, 63, 0, 0, 0),
    @(121, "Ruanda", 60, 0, 0, 60, 0, 0, 0),
    @(122, "Liechtenstein", 56, 0, 0, 56, 0, 0, 0),
    @(123, "Paraguay", 56, 0, 1, 52, 1, 0, 3),
    @(124, "Gibraltar", 56, 0, 14, 42, 0, 0, 0),
    @(125, "Banglades", 48, 0, 15, 28, 1, 0, 5),
    @(126, "Aruba", 46, 0, 1, 45, 0, 0, 0),
    @(127, "Monaco", 43, 0, 1, 41, 0, 0, 1),
    @(128, "Madagascar", 39, 13, 0, 39, 0, 0, 0),
    @(129, "Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
    @(130, "Kenia", 38, 0, 1, 36, 2, 0, 1),
    @(131, "Macao", 34, 0, 10, 24, 0, 0, 0),
    @(132, "Guatemala", 34, 0, 10, 23, 1, 0, 1),
    @(133, "Isla de Man", 32, 0, 0, 32, 0, 0, 0),
    @(134, "Guam", 32, 0, 0, 31, 0, 0, 1),
    @(135, "Jamaica", 32, 2, 2, 29, 0, 0, 1),
    @(136, "Polinesia Francesa", 30, 0, 0, 30, 0, 0, 0),
    @(137, "Uganda", 30, 0, 0, 30, 0, 0, 0),
    @(138, "Zambia", 28, 0, 0, 28, 0, 0, 0),
    @(139, "Guayana Francesa", 28, 0, 6, 22, 0, 0, 0),
    @(140, "Barbados", 26, 0, 0, 26, 0, 0, 0),
    @(141, "Togo", 25, 0, 1, 23, 0, 0, 1),
    @(142, "El Salvador", 24, 5, 0, 24, 0, 0, 0),
    @(143, "Mali", 18, 0, 0, 17, 0, 0, 1),
    @(144, "Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
    @(145, "Bermudas", 17, 0, 2, 15, 0, 0, 0),
    @(146, "Maldivas", 17, 1, 11, 6, 0, 0, 0),
    @(147, "Etiopia", 16, 0, 1, 15, 0, 0, 0),
    @(148, "Nueva Caledonia", 15, 0, 0, 15, 0, 0, 0),
    @(149, "Haiti", 15, 7, 1, 14, 0, 0, 0),
    @(150, "Republica de Yibuti", 14, 0, 0, 14, 0, 0, 0),
    @(151, "Tanzania", 14, 0, 1, 13, 0, 0, 0),
    @(152, "Mongolia", 12, 0, 0, 12, 0, 0, 0),
    @(153, "Guinea Ecuatorial", 12, 0, 0, 12, 0, 0, 0),
    @(154, "Dominica", 11, 0, 0, 11, 0, 0, 0),
    @(155, "San Martin (Parte Francesa)", 11, 0, 0, 11, 0, 0, 0),
    @(156, "Namibia", 11, 3, 2, 9, 0, 0, 0),
    @(157, "Niger", 10, 0, 0, 9, 0, 0, 1),
    @(158, "Bahamas", 10, 0, 1, 9, 0, 0, 0),
    @(159, "Groenlandia", 10, 0, 2, 8, 0, 0, 0),
    @(160, "Suazilandia", 9, 0, 0, 9, 0, 0, 0),
    @(161, "Granada", 9, 2, 0, 9, 0, 0, 0),
    @(162, "Laos", 8, 0, 0, 8, 0, 0, 0),
    @(163, "Seychelles", 8, 0, 0, 8, 0, 0, 0),
    @(164, "Guinea", 8, 0, 0, 8, 0, 0, 0),
    @(165, "Birmania", 8, 0, 0, 8, 0, 0, 0),
    @(166, "Surinam", 8, 0, 0, 8, 0, 0, 0),
    @(167, "Mozambique", 8, 0, 0, 8, 0, 0, 0),
    @(168, "Guyana", 8, 0, 0, 7, 0, 0, 1),
    @(169, "Islas Caimanes", 8, 0, 0, 7, 0, 0, 1),
    @(170, "Curazao", 8, 0, 2, 5, 0, 0, 1),
    @(171, "Antigua y Barbuda", 7, 0, 0, 7, 0, 0, 0),
    @(172, "Gabon", 7, 0, 0, 6, 0, 0, 1),
    @(173, "Zimbabue", 7, 0, 0, 6, 0, 0, 1),
    @(174, "Santa Sede", 6, 0, 0, 6, 0, 0, 0),
    @(175, "Eritrea", 6, 0, 0, 6, 0, 0, 0),
    @(176, "Benin", 6, 0, 0, 6, 0, 0, 0),
    @(177, "Cabo Verde", 6, 0, 0, 5, 0, 0, 1),
    @(178, "Montserrat", 5, 0, 0, 5, 0, 0, 0),
    @(179, "San Bartolome", 5, 0, 0, 5, 0, 0, 0),
    @(180, "Mauritania", 5, 0, 0, 5, 0, 0, 0),
    @(181, "Angola", 5, 0, 0, 5, 0, 0, 0),
    @(182, "Fiyi", 5, 0, 0, 5, 0, 0, 0),
    @(183, "Siria", 5, 0, 0, 5, 0, 0, 0),
    @(184, "Nepal", 5, 0, 1, 4, 0, 0, 0),
    @(185, "Sudan", 5, 0, 0, 4, 0, 0, 1),
    @(186, "Islas Turcas y Caicos", 4, 0, 0, 4, 0, 0, 0),
    @(187, "Butan", 4, 1, 0, 4, 0, 0, 0),
    @(188, "Congo", 4, 0, 0, 4, 0, 0, 0),
    @(189, "Nicaragua", 4, 0, 0, 3, 0, 0, 1),
    @(190, "Republica del Chad", 3, 0, 0, 3, 0, 0, 0),
    @(191, "San Martin (Parte Holandesa)", 3, 0, 0, 3, 0, 0, 0),
    @(192, "Liberia", 3, 0, 0, 3, 0, 0, 0),
    @(193, "Libia", 3, 0, 0, 3, 0, 0, 0),
    @(194, "Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0),
    @(195, "Somalia", 3, 0, 0, 3, 0, 0, 0),
    @(196, "Santa Lucia", 3, 0, 1, 2, 0, 0, 0),
    @(197, "Gambia", 3, 0, 0, 2, 0, 0, 1),
    @(198, "Guinea-Bisau", 2, 0, 0, 2, 0, 0, 0),
    @(199, "San Cristobal y Nieves", 2, 0, 0, 2, 0, 0, 0),
    @(200, "Islas Virgenes Britanicas", 2, 0, 0, 2, 0, 0, 0),
    @(201, "Belice", 2, 0, 0, 2, 0, 0, 0),
    @(202, "Anguila", 2, 0, 0, 2, 0, 0, 0),
    @(203, "Timor Oriental", 1, 0, 0, 1, 0, 0, 0),
    @(204, "Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0),
    @(205, "San Vicente y las Granadinas", 1, 0, 1, 0, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}

"Update complete: " + $data.Count + " rows written"
